$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the current row 143,
# shifting all the subsequent rows (143-206) down by one (to 144-207).
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A143").Value = 8
$ws.Range("B143").Value = "Terminal La Palmera de La Serena"
$ws.Range("C143").Value = "Coquimbo"
$ws.Range("D143").Value = 44510
$ws.Range("E143").Value = 4
$ws.Range("F143").Value = 100112032
$ws.Range("G143").Value = "Zapallo italiano"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 540
$ws.Range("K143").Value = 10000
$ws.Range("L143").Value = 10500
$ws.Range("M143").Value = 10250
$ws.Range("N143").Value = "$/caja 70 unidades"
$ws.Range("O143").Value = "Provincia de Limarí"
$ws.Range("P143").Value = 146
$ws.Range("Q143").Value = 70
$ws.Range("R143").Value = "Hortaliza"
